# Reorder/insert sheets so that the final load order becomes:
#   S2, S0, Some Sheet, S1
# starting from the original order: S1, S2
#
# New sheets:
#   S0         -> A1 = "SX"
#   Some Sheet -> empty

$wb = $excel.ActiveWorkbook

# Insert a new sheet "S0" right after "S2"
$s0 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item("S2"))
$s0.Name = "S0"
$s0.Range("A1").Value = "SX"

# Insert a new, empty sheet "Some Sheet" right after "S0"
$someSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item("S0"))
$someSheet.Name = "Some Sheet"

# Move "S1" to the very end, after "Some Sheet"
$wb.Worksheets.Item("S1").Move($null, $wb.Worksheets.Item("Some Sheet"))

# Make "S0" the active sheet/tab (second tab, index 1 in the saved workbook)
$wb.Worksheets.Item("S0").Activate()
